$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 235 (pushes the existing rows 235:253 down to 236:254)
$ws.Rows("235:235").Insert()

# Populate the newly inserted row 235 with a new weekly price record
# (same market/category/quality/unit/origin/classification as the row that
# used to be at 235, but a new date and a different reported volume)
$ws.Range("A235").Value = 4
$ws.Range("B235").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C235").Value = "Los Lagos"
$ws.Range("D235").Value = 44714
$ws.Range("E235").Value = 10
$ws.Range("F235").Value = 100112017
$ws.Range("G235").Value = "Apio"
$ws.Range("H235").Value = "Americana (o)"
$ws.Range("I235").Value = "Primera"
$ws.Range("J235").Value = 35
$ws.Range("K235").Value = 12000
$ws.Range("L235").Value = 12000
$ws.Range("M235").Value = 12000
$ws.Range("N235").Value = "`$/docena de matas"
$ws.Range("O235").Value = "Región de Coquimbo"
$ws.Range("P235").Value = 2000
$ws.Range("Q235").Value = 6
$ws.Range("R235").Value = "Hortaliza"
